$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 390
$ws1.Range("F5").Value = 2116
$ws1.Range("F6").Value = 10
$ws1.Range("F7").Value = 11013
$ws1.Range("F8").Value = 189
$ws1.Range("F10").Value = 301
$ws1.Range("F12").Value = 10885
$ws1.Range("F14").Value = 1134
$ws1.Range("F15").Value = 26
$ws1.Range("F16").Value = 758
$ws1.Range("F17").Value = 5455
$ws1.Range("F18").Value = 87
$ws1.Range("F19").Value = 3409
$ws1.Range("F20").Value = 5

# Sheet "全部类型" (sheet4) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 390
$ws4.Range("F6").Value = 2116
$ws4.Range("F8").Value = 10
$ws4.Range("F10").Value = 11013
$ws4.Range("F11").Value = 189
$ws4.Range("F13").Value = 301
$ws4.Range("F15").Value = 10885
$ws4.Range("F17").Value = 1134
$ws4.Range("F18").Value = 26
$ws4.Range("F19").Value = 758
$ws4.Range("F20").Value = 5455
$ws4.Range("F21").Value = 87
$ws4.Range("F22").Value = 3409
$ws4.Range("F23").Value = 5
